$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 holds the text "1" (not the number 1) -- a leading apostrophe is the
# normal Excel way to force a numeric-looking entry to be stored as text.
$ws.Range("A2").Value = "'1"
$ws.Range("B2").Value = "Emerging Technology Course-I, Communicative English"
